# 1st May 2nd update
# Inserts a new "01/05/2020" data column into the state-wise COVID-19 long
# table. The new column is inserted before the existing "27/04/2020" column
# (column BI), pushing the last four date columns (27/04, 28/04, 29/04,
# 30/04/2020) one column to the right (BI:BL -> BJ:BM), and the newly freed
# BI column receives the "01/05/2020" figures. The last-reported day
# (30/04/2020, now column BM) is also refreshed with slightly more complete
# totals as extra records came in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at BI - shifts BI:BL (27/04..30/04/2020) right to BJ:BM.
$ws.Columns("BI:BI").Insert()

# New column header. Force text (not an auto-converted date serial) to
# match the plain-text date headers used throughout row 1, then drop back
# to the default "Normal" style so no stray number format sticks around.
$ws.Range("BI1").NumberFormat = "@"
$ws.Range("BI1").Value = "01/05/2020"
$ws.Range("BI1").Style = "Normal"

# New column data (mostly blank/no new cases yet on this first day of reporting).
$ws.Range("BI33").Value = 2

# The trailing column (BM, "30/04/2020") picks up late-arriving case reports,
# so a handful of rows differ from a pure shift of the old BL column.
$ws.Range("BM2").Value = 1
$ws.Range("BM3").Value = 75
$ws.Range("BM5").Value = 2
$ws.Range("BM6").Value = 27
$ws.Range("BM7").Value = 7
$ws.Range("BM8").Value = 4
$ws.Range("BM9").Value = 3
$ws.Range("BM11").Value = 37
$ws.Range("BM12").Value = 9
$ws.Range("BM13").Value = 5
$ws.Range("BM14").Value = 12
$ws.Range("BM15").Value = 1
$ws.Range("BM16").Value = 44
$ws.Range("BM17").Value = 16
$ws.Range("BM19").Value = 16
$ws.Range("BM20").Value = 24
$ws.Range("BM24").Value = 19
$ws.Range("BM26").Value = 18
$ws.Range("BM27").Value = 29
$ws.Range("BM28").Value = 174
$ws.Range("BM29").Value = 5
$ws.Range("BM31").Value = 25
$ws.Range("BM32").Value = 5
$ws.Range("BM33").Value = 3

# Rows 32 and 34 also had their already-shifted 28/04 and 29/04 (BJ/BK)
# figures revised.
$ws.Range("BJ32").Value = ""
$ws.Range("BK32").Value = 2

$ws.Range("BJ34").Value = 2
$ws.Range("BK34").Value = 2
$ws.Range("BM34").Value = ""
